# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Integral"     (the slide master's / every slide's theme)
#   ppt/theme/theme2.xml -> "Office Theme" (the notes master's theme)
#
# The commit swaps the two themes' contents, so the slides now render with
# the "Office Theme" palette (and the notes master would pick up "Integral").
# Font scheme (Arial/Arial) and format scheme (fills/lines/effects) are
# identical between the two themes already, so only the 12 color-scheme
# slots (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) actually need to change.
#
# Apply the "Office Theme" palette to the presentation's design via the
# documented ThemeColorScheme.Colors(index).RGB hook (an OLE_COLOR, packed
# as 0x00BBGGRR).

$p = $ppt.ActivePresentation
$theme = $p.Designs.Item(1).SlideMaster.Theme
$clr = $theme.ThemeColorScheme

$clr.Colors(1).RGB  = 0x000000   # dk1      srgbClr 000000
$clr.Colors(2).RGB  = 0xFFFFFF   # lt1      srgbClr FFFFFF
$clr.Colors(3).RGB  = 0x6A5444   # dk2      srgbClr 44546A
$clr.Colors(4).RGB  = 0xE6E6E7   # lt2      srgbClr E7E6E6
$clr.Colors(5).RGB  = 0xD59B5B   # accent1  srgbClr 5B9BD5
$clr.Colors(6).RGB  = 0x317DED   # accent2  srgbClr ED7D31
$clr.Colors(7).RGB  = 0xA5A5A5   # accent3  srgbClr A5A5A5
$clr.Colors(8).RGB  = 0x00C0FF   # accent4  srgbClr FFC000
$clr.Colors(9).RGB  = 0xC47244   # accent5  srgbClr 4472C4
$clr.Colors(10).RGB = 0x47AD70   # accent6  srgbClr 70AD47
$clr.Colors(11).RGB = 0xC16305   # hlink    srgbClr 0563C1
$clr.Colors(12).RGB = 0x724F95   # folHlink srgbClr 954F72
